$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 34938.21
$ws.Range("I98").Value = 35551.938
$ws.Range("K98").Value = 35551.938
$ws.Range("M98").Value = -34053.938
$ws.Range("H100").Value = 74456.94
$ws.Range("I100").Value = 47485.46
$ws.Range("K100").Value = 47485.46
$ws.Range("M100").Value = -46944.46
$ws.Range("H122").Value = 34938.21
$ws.Range("I122").Value = 35551.938
$ws.Range("K122").Value = 106655.814
$ws.Range("M122").Value = -104205.814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8610.375
$ws.Range("I2").Value = 9501.308000000001
$ws.Range("K2").Value = 9501.308000000001
$ws.Range("M2").Value = -9388.308000000001
$ws.Range("H63").Value = 3450.1667
$ws.Range("I63").Value = 3178
$ws.Range("K63").Value = 3178
$ws.Range("M63").Value = -2492
$ws.Range("H66").Value = 3450.1667
$ws.Range("I66").Value = 3178
$ws.Range("K66").Value = 15890
$ws.Range("M66").Value = -12458
$ws.Range("H116").Value = 8610.375
$ws.Range("I116").Value = 9501.308000000001
$ws.Range("K116").Value = 9501.308000000001
$ws.Range("M116").Value = -7207.308000000001
$ws.Range("H122").Value = 393502.56
$ws.Range("I122").Value = 2802.25
$ws.Range("K122").Value = 8406.75
$ws.Range("M122").Value = -5956.75
$ws.Range("H132").Value = 4877.926
$ws.Range("I132").Value = 4221
$ws.Range("J132").Value = 6191.778
$ws.Range("K132").Value = 12663
$ws.Range("L132").Value = 18575.334
$ws.Range("M132").Value = -10133
$ws.Range("N132").Value = -23635.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8610.375
$ws.Range("I3").Value = 9501.308000000001
$ws.Range("K3").Value = 9501.308000000001
$ws.Range("M3").Value = -9387.308000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 50500500
$ws.Range("I6").Value = 50500500
$ws.Range("K6").Value = 50500500
$ws.Range("M6").Value = -50500387
$ws.Range("H31").Value = 2643.1667
$ws.Range("I31").Value = 1190.8889
$ws.Range("J31").Value = 7000
$ws.Range("K31").Value = 1190.8889
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = -895.8888999999999
$ws.Range("N31").Value = -7590
$ws.Range("H34").Value = 2643.1667
$ws.Range("I34").Value = 1190.8889
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 1190.8889
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = -988.8888999999999
$ws.Range("N34").Value = -7404
$ws.Range("H70").Value = 87499.5
$ws.Range("J70").Value = 87499.5
$ws.Range("L70").Value = 87499.5
$ws.Range("N70").Value = -88129.5
$ws.Range("H73").Value = 87499.5
$ws.Range("J73").Value = 87499.5
$ws.Range("L73").Value = 87499.5
$ws.Range("N73").Value = -89683.5
$ws.Range("H80").Value = 44999
$ws.Range("J80").Value = 44999
$ws.Range("L80").Value = 44999
$ws.Range("N80").Value = -47245
$ws.Range("H83").Value = 44999
$ws.Range("J83").Value = 44999
$ws.Range("L83").Value = 134997
$ws.Range("N83").Value = -146229

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 137.21428
$ws.Range("I33").Value = 106.1
$ws.Range("J33").Value = 215
$ws.Range("K33").Value = 636.5999999999999
$ws.Range("L33").Value = 1290
$ws.Range("M33").Value = -353.5999999999999
$ws.Range("N33").Value = -1856
$ws.Range("H80").Value = 400999.66
$ws.Range("J80").Value = 600000
$ws.Range("L80").Value = 1800000
$ws.Range("N80").Value = -1801872
$ws.Range("H83").Value = 400999.66
$ws.Range("J83").Value = 600000
$ws.Range("L83").Value = 5400000
$ws.Range("N83").Value = -5409360
$ws.Range("H109").Value = 2534.25
$ws.Range("I109").Value = 1897
$ws.Range("K109").Value = 5691
$ws.Range("M109").Value = -4651

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 21399.4
$ws.Range("I99").Value = 19249.5
$ws.Range("K99").Value = 19249.5
$ws.Range("M99").Value = -17003.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3227.2222
$ws.Range("J22").Value = 2508.1428
$ws.Range("L22").Value = 2508.1428
$ws.Range("N22").Value = -3098.1428
$ws.Range("H27").Value = 3227.2222
$ws.Range("J27").Value = 2508.1428
$ws.Range("L27").Value = 2508.1428
$ws.Range("N27").Value = -2722.1428
$ws.Range("H40").Value = 13255.955
$ws.Range("I40").Value = 14174.471
$ws.Range("K40").Value = 14174.471
$ws.Range("M40").Value = -14038.471
$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("M55").Value = -827
$ws.Range("H68").Value = 7143
$ws.Range("I68").Value = 8501
$ws.Range("K68").Value = 8501
$ws.Range("M68").Value = -7752
$ws.Range("H71").Value = 7143
$ws.Range("I71").Value = 8501
$ws.Range("K71").Value = 42505
$ws.Range("M71").Value = -38761
$ws.Range("H74").Value = 70196
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 70196
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 7714.2856
$ws.Range("I100").Value = 30000
$ws.Range("K100").Value = 30000
$ws.Range("M100").Value = -29459

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 172848.1
$ws.Range("I62").Value = 488566.16
$ws.Range("J62").Value = 2846.077
$ws.Range("K62").Value = 488566.16
$ws.Range("L62").Value = 2846.077
$ws.Range("M62").Value = -487942.16
$ws.Range("N62").Value = -4094.077
$ws.Range("H65").Value = 172848.1
$ws.Range("I65").Value = 488566.16
$ws.Range("J65").Value = 2846.077
$ws.Range("K65").Value = 2442830.8
$ws.Range("L65").Value = 14230.385
$ws.Range("M65").Value = -2439710.8
$ws.Range("N65").Value = -20470.385
$ws.Range("H75").Value = 10000
$ws.Range("I75").Value = 10000
$ws.Range("K75").Value = 10000
$ws.Range("M75").Value = -9064
$ws.Range("H78").Value = 10000
$ws.Range("I78").Value = 10000
$ws.Range("K78").Value = 30000
$ws.Range("M78").Value = -25320
$ws.Range("H99").Value = 24999
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 6799.9395
$ws.Range("I122").Value = 4019.0476
$ws.Range("K122").Value = 12057.1428
$ws.Range("M122").Value = -9607.1428
$ws.Range("H132").Value = 20166.637
$ws.Range("I132").Value = 25518.61
$ws.Range("J132").Value = 7857.1
$ws.Range("K132").Value = 76555.83
$ws.Range("L132").Value = 23571.3
$ws.Range("M132").Value = -74025.83
$ws.Range("N132").Value = -28631.3
$ws.Range("H133").Value = 47999.668
$ws.Range("J133").Value = 47999.668
$ws.Range("L133").Value = 47999.668
$ws.Range("N133").Value = -58119.668
